# Auto-assisted generation: applies the Halicarnassus_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Cells.Item(11, 8).Value = 43   # H11: 41.333332 -> 43
$ws.Cells.Item(11, 9).Value = 43   # I11: 41.333332 -> 43
$ws.Cells.Item(11, 11).Value = 43   # K11: 41.333332 -> 43
$ws.Cells.Item(11, 13).Value = 97   # M11: 98.666668 -> 97

# Row 87
$ws.Cells.Item(87, 8).Value = 57855   # H87: 38284 -> 57855
$ws.Cells.Item(87, 9).Value = 12000   # I87: 10666.667 -> 12000
$ws.Cells.Item(87, 10).Value = 73140   # J87: 79710 -> 73140
$ws.Cells.Item(87, 11).Value = 12000   # K87: 10666.667 -> 12000
$ws.Cells.Item(87, 12).Value = 73140   # L87: 79710 -> 73140
$ws.Cells.Item(87, 13).Value = -10752   # M87: -9418.666999999999 -> -10752
$ws.Cells.Item(87, 14).Value = -75636   # N87: -82206 -> -75636

# Row 90
$ws.Cells.Item(90, 8).Value = 57855   # H90: 38284 -> 57855
$ws.Cells.Item(90, 9).Value = 12000   # I90: 10666.667 -> 12000
$ws.Cells.Item(90, 10).Value = 73140   # J90: 79710 -> 73140
$ws.Cells.Item(90, 11).Value = 36000   # K90: 32000.001 -> 36000
$ws.Cells.Item(90, 12).Value = 219420   # L90: 239130 -> 219420
$ws.Cells.Item(90, 13).Value = -29760   # M90: -25760.001 -> -29760
$ws.Cells.Item(90, 14).Value = -231900   # N90: -251610 -> -231900

# Row 111
$ws.Cells.Item(111, 8).Value = 1827.5714   # H111: 1762 -> 1827.5714
$ws.Cells.Item(111, 9).Value = 1827.5714   # I111: 1756.5714 -> 1827.5714
$ws.Cells.Item(111, 10).Value = 0   # J111: 1800 -> 0
$ws.Cells.Item(111, 11).Value = 5482.7142   # K111: 5269.7142 -> 5482.7142
$ws.Cells.Item(111, 12).Value = 0   # L111: 5400 -> 0
$ws.Cells.Item(111, 13).Value = None   # M111: -2202.7142 -> None
$ws.Cells.Item(111, 14).ClearContents()   # N111: remove (was -11534)

# Row 132
$ws.Cells.Item(132, 8).Value = 20183.154   # H132: 20615.166 -> 20183.154
$ws.Cells.Item(132, 10).Value = 11398.5   # J132: 10198.333 -> 11398.5
$ws.Cells.Item(132, 12).Value = 34195.5   # L132: 30594.999 -> 34195.5
$ws.Cells.Item(132, 14).Value = -39255.5   # N132: -35654.999 -> -39255.5

# Row 137
$ws.Cells.Item(137, 8).Value = 3552.5557   # H137: 3526.2942 -> 3552.5557
$ws.Cells.Item(137, 10).Value = 4989.1113   # J137: 5112.875 -> 4989.1113
$ws.Cells.Item(137, 12).Value = 14967.3339   # L137: 15338.625 -> 14967.3339
$ws.Cells.Item(137, 14).Value = -20067.3339   # N137: -20438.625 -> -20067.3339

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 2958.4   # H2: 2705 -> 2958.4
$ws.Cells.Item(2, 9).Value = 2237.6   # I2: 2052.2727 -> 2237.6
$ws.Cells.Item(2, 10).Value = 4400   # J2: 4500 -> 4400
$ws.Cells.Item(2, 11).Value = 2237.6   # K2: 2052.2727 -> 2237.6
$ws.Cells.Item(2, 12).Value = 4400   # L2: 4500 -> 4400
$ws.Cells.Item(2, 13).Value = -2124.6   # M2: -1939.2727 -> -2124.6
$ws.Cells.Item(2, 14).Value = -4626   # N2: -4726 -> -4626

# Row 61
$ws.Cells.Item(61, 8).Value = 5562.727   # H61: 5639.9 -> 5562.727
$ws.Cells.Item(61, 9).Value = 4027.1428   # I61: 4199.857 -> 4027.1428
$ws.Cells.Item(61, 10).Value = 8250   # J61: 9000 -> 8250
$ws.Cells.Item(61, 11).Value = 4027.1428   # K61: 4199.857 -> 4027.1428
$ws.Cells.Item(61, 12).Value = 8250   # L61: 9000 -> 8250
$ws.Cells.Item(61, 13).Value = -3815.1428   # M61: -3987.857 -> -3815.1428
$ws.Cells.Item(61, 14).Value = -8674   # N61: -9424 -> -8674

# Row 74
$ws.Cells.Item(74, 8).Value = 2873   # H74: 2998 -> 2873
$ws.Cells.Item(74, 9).Value = 2832   # I74: 2998 -> 2832
$ws.Cells.Item(74, 10).Value = 2996   # J74: 0 -> 2996
$ws.Cells.Item(74, 11).Value = 2832   # K74: 2998 -> 2832
$ws.Cells.Item(74, 12).Value = 2996   # L74: 0 -> 2996
$ws.Cells.Item(74, 13).Value = -1958   # M74: -2124 -> -1958
$ws.Cells.Item(74, 14).Value = -4744   # N74: None -> -4744

# Row 77
$ws.Cells.Item(77, 8).Value = 2873   # H77: 2998 -> 2873
$ws.Cells.Item(77, 9).Value = 2832   # I77: 2998 -> 2832
$ws.Cells.Item(77, 10).Value = 2996   # J77: 0 -> 2996
$ws.Cells.Item(77, 11).Value = 14160   # K77: 14990 -> 14160
$ws.Cells.Item(77, 12).Value = 14980   # L77: 0 -> 14980
$ws.Cells.Item(77, 13).Value = -9792   # M77: -10622 -> -9792
$ws.Cells.Item(77, 14).Value = -23716   # N77: None -> -23716

# Row 116
$ws.Cells.Item(116, 8).Value = 2958.4   # H116: 2705 -> 2958.4
$ws.Cells.Item(116, 9).Value = 2237.6   # I116: 2052.2727 -> 2237.6
$ws.Cells.Item(116, 10).Value = 4400   # J116: 4500 -> 4400
$ws.Cells.Item(116, 11).Value = 2237.6   # K116: 2052.2727 -> 2237.6
$ws.Cells.Item(116, 12).Value = 4400   # L116: 4500 -> 4400
$ws.Cells.Item(116, 13).Value = 56.40000000000009   # M116: 241.7273 -> 56.40000000000009
$ws.Cells.Item(116, 14).Value = -8988   # N116: -9088 -> -8988

# Row 132
$ws.Cells.Item(132, 8).Value = 3209.8462   # H132: 3403.7083 -> 3209.8462
$ws.Cells.Item(132, 9).Value = 2585.0435   # I132: 2747.0952 -> 2585.0435
$ws.Cells.Item(132, 11).Value = 7755.130500000001   # K132: 8241.285600000001 -> 7755.130500000001
$ws.Cells.Item(132, 13).Value = -5225.130500000001   # M132: -5711.285600000001 -> -5225.130500000001

# Row 136
$ws.Cells.Item(136, 8).Value = 5562.727   # H136: 5639.9 -> 5562.727
$ws.Cells.Item(136, 9).Value = 4027.1428   # I136: 4199.857 -> 4027.1428
$ws.Cells.Item(136, 10).Value = 8250   # J136: 9000 -> 8250
$ws.Cells.Item(136, 11).Value = 12081.4284   # K136: 12599.571 -> 12081.4284
$ws.Cells.Item(136, 12).Value = 24750   # L136: 27000 -> 24750
$ws.Cells.Item(136, 13).Value = -9531.428400000001   # M136: -10049.571 -> -9531.428400000001
$ws.Cells.Item(136, 14).Value = -29850   # N136: -32100 -> -29850

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 2958.4   # H3: 2705 -> 2958.4
$ws.Cells.Item(3, 9).Value = 2237.6   # I3: 2052.2727 -> 2237.6
$ws.Cells.Item(3, 10).Value = 4400   # J3: 4500 -> 4400
$ws.Cells.Item(3, 11).Value = 2237.6   # K3: 2052.2727 -> 2237.6
$ws.Cells.Item(3, 12).Value = 4400   # L3: 4500 -> 4400
$ws.Cells.Item(3, 13).Value = -2123.6   # M3: -1938.2727 -> -2123.6
$ws.Cells.Item(3, 14).Value = -4628   # N3: -4728 -> -4628

# Row 6
$ws.Cells.Item(6, 8).Value = 0   # H6: 27500 -> 0
$ws.Cells.Item(6, 10).Value = 0   # J6: 27500 -> 0
$ws.Cells.Item(6, 12).Value = None   # L6: 27500 -> None
$ws.Cells.Item(6, 14).ClearContents()   # N6: remove (was -27726)

# Row 57
$ws.Cells.Item(57, 8).Value = 94999   # H57: 0 -> 94999
$ws.Cells.Item(57, 10).Value = 94999   # J57: 0 -> 94999
$ws.Cells.Item(57, 12).Value = 94999   # L57: 0 -> 94999
$ws.Cells.Item(57, 14).Value = -96439   # N57: None -> -96439

# Row 58
$ws.Cells.Item(58, 8).Value = 150000   # H58: 63853.332 -> 150000
$ws.Cells.Item(58, 10).Value = 150000   # J58: 63853.332 -> 150000
$ws.Cells.Item(58, 12).Value = 150000   # L58: 63853.332 -> 150000
$ws.Cells.Item(58, 14).Value = -150588   # N58: -64441.332 -> -150588

# Row 59
$ws.Cells.Item(59, 8).Value = 94999   # H59: 0 -> 94999
$ws.Cells.Item(59, 10).Value = 94999   # J59: 0 -> 94999
$ws.Cells.Item(59, 12).Value = 94999   # L59: 0 -> 94999
$ws.Cells.Item(59, 14).Value = -96693   # N59: None -> -96693

# Row 80
$ws.Cells.Item(80, 8).Value = 213.22223   # H80: 219 -> 213.22223
$ws.Cells.Item(80, 9).Value = 173.16667   # I80: 184.8 -> 173.16667
$ws.Cells.Item(80, 11).Value = 173.16667   # K80: 184.8 -> 173.16667
$ws.Cells.Item(80, 13).Value = 824.8333299999999   # M80: 813.2 -> 824.8333299999999

# Row 83
$ws.Cells.Item(83, 8).Value = 213.22223   # H83: 219 -> 213.22223
$ws.Cells.Item(83, 9).Value = 173.16667   # I83: 184.8 -> 173.16667
$ws.Cells.Item(83, 11).Value = 865.8333500000001   # K83: 924 -> 865.8333500000001
$ws.Cells.Item(83, 13).Value = 4126.16665   # M83: 4068 -> 4126.16665

# Row 105
$ws.Cells.Item(105, 8).Value = 1570.8   # H105: 2296.6667 -> 1570.8
$ws.Cells.Item(105, 9).Value = 1570.8   # I105: 2296.6667 -> 1570.8
$ws.Cells.Item(105, 11).Value = 1570.8   # K105: 2296.6667 -> 1570.8
$ws.Cells.Item(105, 13).Value = 176.2   # M105: -549.6667000000002 -> 176.2

# Row 134
$ws.Cells.Item(134, 8).Value = 7189.1   # H134: 8643.5 -> 7189.1
$ws.Cells.Item(134, 9).Value = 6876.8887   # I134: 8372.4 -> 6876.8887
$ws.Cells.Item(134, 11).Value = 20630.6661   # K134: 25117.2 -> 20630.6661
$ws.Cells.Item(134, 13).Value = -18095.6661   # M134: -22582.2 -> -18095.6661

# Row 136
$ws.Cells.Item(136, 8).Value = 94999   # H136: 0 -> 94999
$ws.Cells.Item(136, 10).Value = 94999   # J136: 0 -> 94999
$ws.Cells.Item(136, 12).Value = 94999   # L136: 0 -> 94999
$ws.Cells.Item(136, 14).Value = -105199   # N136: None -> -105199

# Row 139
$ws.Cells.Item(139, 8).Value = 84999   # H139: 0 -> 84999
$ws.Cells.Item(139, 10).Value = 84999   # J139: 0 -> 84999
$ws.Cells.Item(139, 12).Value = 84999   # L139: 0 -> 84999
$ws.Cells.Item(139, 14).Value = -95279   # N139: None -> -95279

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Cells.Item(58, 8).Value = 3522.5557   # H58: 2880.1177 -> 3522.5557
$ws.Cells.Item(58, 10).Value = 14444   # J58: 0 -> 14444
$ws.Cells.Item(58, 12).Value = 14444   # L58: 0 -> 14444
$ws.Cells.Item(58, 14).Value = -14850   # N58: None -> -14850

# Row 99
$ws.Cells.Item(99, 8).Value = 3839.8   # H99: 3910.3684 -> 3839.8
$ws.Cells.Item(99, 9).Value = 3737.3125   # I99: 3819.8667 -> 3737.3125
$ws.Cells.Item(99, 11).Value = 3737.3125   # K99: 3819.8667 -> 3737.3125
$ws.Cells.Item(99, 13).Value = -2239.3125   # M99: -2321.8667 -> -2239.3125

# Row 126
$ws.Cells.Item(126, 8).Value = 3839.8   # H126: 3910.3684 -> 3839.8
$ws.Cells.Item(126, 9).Value = 3737.3125   # I126: 3819.8667 -> 3737.3125
$ws.Cells.Item(126, 11).Value = 11211.9375   # K126: 11459.6001 -> 11211.9375
$ws.Cells.Item(126, 13).Value = -8741.9375   # M126: -8989.6001 -> -8741.9375

# Row 134
$ws.Cells.Item(134, 8).Value = 2538.125   # H134: 5635.4 -> 2538.125
$ws.Cells.Item(134, 9).Value = 2310.348   # I134: 4450 -> 2310.348
$ws.Cells.Item(134, 10).Value = 7777   # J134: 6425.6665 -> 7777
$ws.Cells.Item(134, 11).Value = 6931.044   # K134: 13350 -> 6931.044
$ws.Cells.Item(134, 12).Value = 23331   # L134: 19276.9995 -> 23331
$ws.Cells.Item(134, 13).Value = -4396.044   # M134: -10815 -> -4396.044
$ws.Cells.Item(134, 14).Value = -28401   # N134: -24346.9995 -> -28401

# Row 136
$ws.Cells.Item(136, 8).Value = 3522.5557   # H136: 2880.1177 -> 3522.5557
$ws.Cells.Item(136, 10).Value = 14444   # J136: 0 -> 14444
$ws.Cells.Item(136, 12).Value = 43332   # L136: 0 -> 43332
$ws.Cells.Item(136, 14).Value = -48432   # N136: None -> -48432

$ws = $wb.Worksheets.Item("GSM")
# Row 96
$ws.Cells.Item(96, 8).Value = 46130.5   # H96: 33630.5 -> 46130.5
$ws.Cells.Item(96, 10).Value = 46130.5   # J96: 33630.5 -> 46130.5
$ws.Cells.Item(96, 12).Value = 46130.5   # L96: 33630.5 -> 46130.5
$ws.Cells.Item(96, 14).Value = -51622.5   # N96: -39122.5 -> -51622.5

# Row 132
$ws.Cells.Item(132, 8).Value = 45103.37   # H132: 48411.76 -> 45103.37
$ws.Cells.Item(132, 9).Value = 51436.87   # I132: 55978.617 -> 51436.87
$ws.Cells.Item(132, 11).Value = 154310.61   # K132: 167935.851 -> 154310.61
$ws.Cells.Item(132, 13).Value = -151780.61   # M132: -165405.851 -> -151780.61

# Row 135
$ws.Cells.Item(135, 8).Value = 230000   # H135: 202695 -> 230000
$ws.Cells.Item(135, 10).Value = 230000   # J135: 202695 -> 230000
$ws.Cells.Item(135, 12).Value = 230000   # L135: 202695 -> 230000
$ws.Cells.Item(135, 14).Value = -240140   # N135: -212835 -> -240140

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 629.6   # H16: 699.3333 -> 629.6
$ws.Cells.Item(16, 9).Value = 328   # I16: 359 -> 328
$ws.Cells.Item(16, 10).Value = 1333.3334   # J16: 1124.75 -> 1333.3334
$ws.Cells.Item(16, 11).Value = 328   # K16: 359 -> 328
$ws.Cells.Item(16, 12).Value = 1333.3334   # L16: 1124.75 -> 1333.3334
$ws.Cells.Item(16, 13).Value = -158   # M16: -189 -> -158
$ws.Cells.Item(16, 14).Value = -1673.3334   # N16: -1464.75 -> -1673.3334

# Row 38
$ws.Cells.Item(38, 8).Value = 32999   # H38: 30343.666 -> 32999
$ws.Cells.Item(38, 10).Value = 32999   # J38: 30343.666 -> 32999
$ws.Cells.Item(38, 12).Value = 32999   # L38: 30343.666 -> 32999
$ws.Cells.Item(38, 14).Value = -33819   # N38: -31163.666 -> -33819

# Row 132
$ws.Cells.Item(132, 8).Value = 3891.7036   # H132: 3906.1035 -> 3891.7036
$ws.Cells.Item(132, 9).Value = 2677   # I132: 2762.25 -> 2677
$ws.Cells.Item(132, 10).Value = 4499.0557   # J132: 4341.857 -> 4499.0557
$ws.Cells.Item(132, 11).Value = 8031   # K132: 8286.75 -> 8031
$ws.Cells.Item(132, 12).Value = 13497.1671   # L132: 13025.571 -> 13497.1671
$ws.Cells.Item(132, 13).Value = -5501   # M132: -5756.75 -> -5501
$ws.Cells.Item(132, 14).Value = -18557.1671   # N132: -18085.571 -> -18557.1671

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 1000   # H81: 1066.6666 -> 1000

# Row 84
$ws.Cells.Item(84, 8).Value = 1000   # H84: 1066.6666 -> 1000

# Row 95
$ws.Cells.Item(95, 8).Value = 33833.168   # H95: 31857 -> 33833.168
$ws.Cells.Item(95, 10).Value = 33833.168   # J95: 31857 -> 33833.168
$ws.Cells.Item(95, 12).Value = 33833.168   # L95: 31857 -> 33833.168
$ws.Cells.Item(95, 14).Value = -39325.168   # N95: -37349 -> -39325.168
